$d = $word.ActiveDocument

# wdHeaderFooterIndex constants
$wdHeaderFooterPrimary   = 1
$wdHeaderFooterFirstPage = 2

$sec = $d.Sections.Item(1)

# --- Primary footer: Pearson Edexcel logo (docPr id="1") ---
$footerPrimary = $sec.Footers.Item($wdHeaderFooterPrimary)
if ($footerPrimary.Exists -and $footerPrimary.Range.InlineShapes.Count -ge 1) {
    $shp = $footerPrimary.Range.InlineShapes.Item(1)
    if ($shp.AlternativeText -eq "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png") {
        $shp.Name = "image2.png"
    }
}

# --- First-page footer: Pearson Edexcel logo (docPr id="2") ---
$footerFirst = $sec.Footers.Item($wdHeaderFooterFirstPage)
if ($footerFirst.Exists -and $footerFirst.Range.InlineShapes.Count -ge 1) {
    $shp = $footerFirst.Range.InlineShapes.Item(1)
    if ($shp.AlternativeText -eq "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png") {
        $shp.Name = "image2.png"
    }
}

# --- First-page header: BTEC logo (docPr id="3") ---
$headerFirst = $sec.Headers.Item($wdHeaderFooterFirstPage)
if ($headerFirst.Exists -and $headerFirst.Range.InlineShapes.Count -ge 1) {
    $shp = $headerFirst.Range.InlineShapes.Item(1)
    if ($shp.AlternativeText -eq "BTec_Logo-Orange") {
        $shp.Name = "image1.jpg"
    }
}
